$wb = $excel.ActiveWorkbook

# --- Overview sheet (new row for 6ab923b6-2220-4837-9085-28359a473124.md) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A8").Value = "6ab923b6-2220-4837-9085-28359a473124.md"
$wsOverview.Range("B8").Value = "e2e\6ab923b6-2220-4837-9085-28359a473124.md"
$wsOverview.Range("C8").Value = ".md"
$wsOverview.Range("D8").Value = ""
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = "2016-10-24 09:09:42"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B8"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab923b6-placeholder/e2e/6ab923b6-2220-4837-9085-28359a473124.md",
    "",
    "",
    "e2e\6ab923b6-2220-4837-9085-28359a473124.md"
) | Out-Null

# --- zh-cn sheet (new row for 6ab923b6-2220-4837-9085-28359a473124.md) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A8").Value = "6ab923b6-2220-4837-9085-28359a473124.md"
$wsZh.Range("B8").Value = ".md"
$wsZh.Range("C8").Value = "Ready for handoff"
$wsZh.Range("D8").Value = "e2e"
$wsZh.Range("E8").Value = "ht"
$wsZh.Range("F8").Value = "False"
$wsZh.Range("G8").Value = "6ab923b6-2220-4837-9085-28359a473124.86832aa69635ebbc4d8609dc9f9980bd5e3f700b.zh-cn.xlf"
$wsZh.Range("H8").Value = "2016-10-24 09:09:30"
$wsZh.Range("I8").Value = ""
$wsZh.Range("J8").Value = ""
$wsZh.Range("K8").Value = "0001-01-01 00:00:00"
$wsZh.Range("L8").Value = ""
$wsZh.Range("M8").Value = "True"
$wsZh.Range("N8").Value = ""
$wsZh.Range("O8").Value = "False"
$wsZh.Range("P8").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A8"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab923b6-placeholder/e2e/6ab923b6-2220-4837-9085-28359a473124.md",
    "",
    "",
    "6ab923b6-2220-4837-9085-28359a473124.md"
) | Out-Null

# --- de-de sheet (new row for 6ab923b6-2220-4837-9085-28359a473124.md) ---
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A8").Value = "6ab923b6-2220-4837-9085-28359a473124.md"
$wsDe.Range("B8").Value = ".md"
$wsDe.Range("C8").Value = "Ready for handoff"
$wsDe.Range("D8").Value = "e2e"
$wsDe.Range("E8").Value = "ht"
$wsDe.Range("F8").Value = "False"
$wsDe.Range("G8").Value = "6ab923b6-2220-4837-9085-28359a473124.86832aa69635ebbc4d8609dc9f9980bd5e3f700b.de-de.xlf"
$wsDe.Range("H8").Value = "2016-10-24 09:09:42"
$wsDe.Range("I8").Value = ""
$wsDe.Range("J8").Value = ""
$wsDe.Range("K8").Value = "0001-01-01 00:00:00"
$wsDe.Range("L8").Value = ""
$wsDe.Range("M8").Value = "True"
$wsDe.Range("N8").Value = ""
$wsDe.Range("O8").Value = "False"
$wsDe.Range("P8").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A8"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab923b6-placeholder/e2e/6ab923b6-2220-4837-9085-28359a473124.md",
    "",
    "",
    "6ab923b6-2220-4837-9085-28359a473124.md"
) | Out-Null
